# API List.xlsx - "optimaze the search api"
# Adds a new row (row 6) describing the "GET /api/skills/{member id}" endpoint.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$jsonText = @"
    {
        "id": 1, 
        "name": "API Implementation", 
        "classification": "", 
        "prerequisites": null, 
        "knowledge_area": "Software Construction", 
        "rationale": null, 
        "roles_for_skill": null, 
        "related_activities": null, 
        "real_world_scenario": null, 
        "role_of_academia": null, 
        "tools": null, 
        "self_assessment": null, 
        "reference": null, 
        "student_name": "Baiyu Huo", 
        "student_no": 40076004
    }, 
"@

$ws.Range("A6").Value = "A05"
$ws.Range("B6").Value = "/api/skills/{member id}"
$ws.Range("C6").Value = "GET"
$ws.Range("D6").Value = "id"
$ws.Range("E6").Value = $jsonText

$ws.Rows.Item(6).RowHeight = 234.6

$ws.Range("F4").Select()

$wb.Windows.Item(1).WindowState = -4140
